# Actualización desde MV -datos-
# Append the new quarterly data point (01-07-2021) as row 76 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 76

# Column A holds a date-like label ("dd-mm-aaaa") that must stay plain text,
# matching every other cell in the column. Force the cell to Text first so
# the engine doesn't silently convert "01-07-2021" into a date serial, then
# drop the formatting back to Normal so the cell ends up styled exactly like
# its neighbours (A2:A75 carry no explicit style).
$ws.Cells.Item($newRow, 1).NumberFormat = "@"
$ws.Cells.Item($newRow, 1).Value = "01-07-2021"
$ws.Cells.Item($newRow, 1).Style = "Normal"

$ws.Cells.Item($newRow, 2).Value = 3604
$ws.Cells.Item($newRow, 3).Value = 0
$ws.Cells.Item($newRow, 4).Value = 0
$ws.Cells.Item($newRow, 5).Value = 0
$ws.Cells.Item($newRow, 6).Value = 3604
$ws.Cells.Item($newRow, 7).Value = 10
$ws.Cells.Item($newRow, 8).Value = 87
$ws.Cells.Item($newRow, 9).Value = 3506
